$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update A2:A11 test_number values from 2 to 4
$ws.Range("A2:A11").Value = 4

# Update the selected/active cell to J10 (matches sheetView selection in diff)
$ws.Range("J10").Select()
